$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) assignments to be stored as literal text,
# matching the source data (t="inlineStr") rather than being smart-
# parsed into numbers by Excel (e.g. "44.284.29" or "0.678").

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.284.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.363.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.678'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.49'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +12.76%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.546'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +20.02%  '
$ws.Range('E10').Value = '  +7.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.51'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.07%  '
$ws.Range('E12').Value = '  +3.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.714.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '16.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.72'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.904'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.364.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.120.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.27%  '
$ws.Range('E19').Value = '  +5.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.10'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.25%  '
$ws.Range('E21').Value = '  +4.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '256.21'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.84%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.89%  '
$ws.Range('E25').Value = '  +3.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.83%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  +3.68%  '
$ws.Range('E32').Value = '  +5.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.22'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0738'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.17%  '
$ws.Range('E36').Value = '  +10.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.45'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.49'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  +7.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.86%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0982'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.28%  '
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('E48').Value = '  +13.09%  '
$ws.Range('E49').Value = '  +5.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.442.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.53%  '
